$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 496, shifting existing rows 496:615 down to 497:616
$ws.Rows(496).Insert()

# Populate the newly inserted row 496 with its data
$ws.Range("A496").Value = 10
$ws.Range("B496").Value = "Vega Modelo de Temuco"
$ws.Range("C496").Value = "La Araucanía"
$ws.Range("D496").Value = 44995
$ws.Range("E496").Value = 9
$ws.Range("F496").Value = 100112023
$ws.Range("G496").Value = "Brócoli"
$ws.Range("H496").Value = "Sin especificar"
$ws.Range("I496").Value = "Primera"
$ws.Range("J496").Value = 650
$ws.Range("K496").Value = 1300
$ws.Range("L496").Value = 1300
$ws.Range("M496").Value = 1300
$ws.Range("N496").Value = "$/unidad"
$ws.Range("O496").Value = "Provincia de Cautín"
$ws.Range("P496").Value = 1300
$ws.Range("Q496").Value = 1
$ws.Range("R496").Value = "Hortaliza"
